$d = $word.ActiveDocument

function Find-And-Highlight([string]$text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.HighlightColorIndex = 7
    }
    return $found
}

# --- Highlight the grading-criteria paragraphs (yellow) ---
Find-And-Highlight("Asynchrone (affichage des cartes, fin de partie) : 10 %")
Find-And-Highlight("Validations : 20 %")
Find-And-Highlight("Jeu fonctionnel (gestion des clics, fin de partie, etc.): 30 %")
Find-And-Highlight("Test fonctionnel : 10 %")
Find-And-Highlight("Utilisation d" + [char]8217 + "un objet : 10 %")
Find-And-Highlight("Utilisation de git. Au moins 3 commit par " + [char]233 + "tudiant: 10 %")
Find-And-Highlight("Visuel agr" + [char]233 + "able : bonus maximal de 10 %")

Write-Output "highlight done"
